$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the player roster rows (A2:C18) to match the updated sheet.
# Row 1 (headers) and several player rows stay the same; the rows for
# Damian Lillard, LaMelo Ball and Tari Eason move up right after
# Luke Kennard, and Devin Vassell moves down after Collin Sexton.

$ws.Range("A2").Value = "Derrick White"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Boston Celtics"

$ws.Range("A3").Value = "Luke Kennard"
$ws.Range("B3").Value = "SG"
$ws.Range("C3").Value = "Memphis Grizzlies"

$ws.Range("A4").Value = "Damian Lillard"
$ws.Range("B4").Value = "PG"
$ws.Range("C4").Value = "Milwaukee Bucks"

$ws.Range("A5").Value = "LaMelo Ball"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Charlotte Hornets"

$ws.Range("A6").Value = "Tari Eason"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Houston Rockets"

$ws.Range("A7").Value = "Onyeka Okongwu"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Atlanta Hawks"

$ws.Range("A8").Value = "Malik Monk"
$ws.Range("B8").Value = "PG,SG,SF"
$ws.Range("C8").Value = "Sacramento Kings"

$ws.Range("A9").Value = "Naz Reid"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Minnesota Timberwolves"

$ws.Range("A10").Value = "Anthony Davis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Los Angeles Lakers"

$ws.Range("A11").Value = "Julius Randle"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Minnesota Timberwolves"

$ws.Range("A12").Value = "Cade Cunningham"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Detroit Pistons"

$ws.Range("A13").Value = "Collin Sexton"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Utah Jazz"

$ws.Range("A14").Value = "Devin Vassell"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "San Antonio Spurs"

$ws.Range("A15").Value = "Isaiah Hartenstein"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Oklahoma City Thunder"

$ws.Range("A16").Value = "Deandre Ayton"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Portland Trail Blazers"

$ws.Range("A17").Value = "Coby White"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Chicago Bulls"

$ws.Range("A18").Value = "Cameron Johnson"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Brooklyn Nets"
